# Scheduled runner update: refresh Market Board pricing / profit figures
# across the Tiamat_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR). Each sheet has the same header layout:
#   A Leve Name | B Leve Item | C Leve Level | D Leve EXP | E Leve Gil
#   F Leve Amount | G Leve Item ID | H currentAveragePrice
#   I currentAveragePriceNQ | J currentAveragePriceHQ | K LevePriceNQ
#   L LevePriceHQ | M LeveProfitNQ | N LeveProfitHQ

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# ---------------- ALC ----------------
Set-Cells "ALC" 15 @{
    H = 1527.9231; I = 1527.9231; K = 4583.7693; M = -4414.7693
}
Set-Cells "ALC" 55 @{
    H = 2003.6666; J = 3117.1667; L = 3117.1667; N = -3545.1667
}
Set-Cells "ALC" 103 @{
    H = 326; I = 326; J = 0; K = 978; L = 0; M = -392
}
(($wb.Worksheets.Item("ALC")).Range("N103")).Value = ""
Set-Cells "ALC" 111 @{
    H = 10001; I = 10001; K = 30003; M = -26936
}
Set-Cells "ALC" 112 @{
    H = 1329.59; J = 1351.6562; L = 4054.9686; N = -6270.9686
}

# ---------------- ARM ----------------
Set-Cells "ARM" 32 @{
    H = 140674.28; I = 141968.31; J = 129027.875; K = 141968.31
    L = 129027.875; M = -141681.31; N = -129601.875
}
Set-Cells "ARM" 45 @{
    H = 3156.7778; I = 2635.1667; K = 2635.1667; M = -2258.1667
}

# ---------------- BSM ----------------
Set-Cells "BSM" 94 @{
    H = 1878.25; I = 1139.875; J = 3355; K = 1139.875; L = 3355
    M = -688.875; N = -4257
}

# ---------------- CRP ----------------
Set-Cells "CRP" 16 @{
    H = 1200; I = 0; J = 1200; K = 0; L = 1200; N = -1774
}
(($wb.Worksheets.Item("CRP")).Range("M16")).Value = ""
Set-Cells "CRP" 39 @{
    H = 22834; I = 2135; J = 36633.332; K = 2135; L = 36633.332
    M = -1744; N = -37415.332
}
Set-Cells "CRP" 49 @{
    H = 22834; I = 2135; J = 36633.332; K = 2135; L = 36633.332
    M = -1953; N = -36997.332
}
Set-Cells "CRP" 105 @{
    H = 719.25; I = 683; K = 683; M = 1064
}
Set-Cells "CRP" 113 @{
    H = 1200; I = 0; J = 1200; K = 0; L = 1200; N = -5540
}
(($wb.Worksheets.Item("CRP")).Range("M113")).Value = ""

# ---------------- CUL ----------------
Set-Cells "CUL" 60 @{
    H = 373.6; I = 373.6; K = 1120.8; M = -869.8000000000002
}
Set-Cells "CUL" 86 @{
    H = 507; I = 414; J = 600; K = 1242; L = 1800; M = -56; N = -4172
}
Set-Cells "CUL" 87 @{
    H = 63948.05; I = 3500; J = 71059.586; K = 10500; L = 213178.758
    M = -9252; N = -215674.758
}
Set-Cells "CUL" 89 @{
    H = 507; I = 414; J = 600; K = 3726; L = 5400; M = 2202; N = -17256
}
Set-Cells "CUL" 90 @{
    H = 63948.05; I = 3500; J = 71059.586; K = 31500; L = 639536.274
    M = -25260; N = -652016.274
}

# ---------------- GSM ----------------
Set-Cells "GSM" 34 @{
    H = 11200; J = 11200; L = 11200; N = -11736
}
Set-Cells "GSM" 76 @{
    H = 11200; J = 11200; L = 11200; N = -11830
}
Set-Cells "GSM" 79 @{
    H = 11200; J = 11200; L = 11200; N = -13384
}
Set-Cells "GSM" 107 @{
    H = 668.7857; I = 766.6667; J = 642.0909; K = 766.6667; L = 642.0909
    M = 1153.3333; N = -4482.0909
}
Set-Cells "GSM" 113 @{
    H = 1400; I = 1400; K = 1400; M = 770
}

# ---------------- LTW ----------------
Set-Cells "LTW" 55 @{
    H = 578.0789; I = 308.29166; K = 308.29166; M = -135.29166
}
Set-Cells "LTW" 68 @{
    H = 2022.3636; I = 1639.2858; J = 2692.75; K = 1639.2858; L = 2692.75
    M = -890.2858000000001; N = -4190.75
}
Set-Cells "LTW" 71 @{
    H = 2022.3636; I = 1639.2858; J = 2692.75; K = 8196.429; L = 13463.75
    M = -4452.429; N = -20951.75
}

# ---------------- WVR ----------------
Set-Cells "WVR" 81 @{
    H = 1497.2222; I = 821.6; J = 1757.0769; K = 1643.2; L = 3514.1538
    M = -582.2; N = -5636.1538
}
Set-Cells "WVR" 84 @{
    H = 1497.2222; I = 821.6; J = 1757.0769; K = 8216; L = 17570.769
    M = -2912; N = -28178.769
}
